$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $origStyle = $Cell.Style
    $Cell.Value = "'" + $Text
    $Cell.Style = $origStyle
}

# Row 2
$ws.Cells.Item(2, 4).Value = "61.030.91"
$ws.Cells.Item(2, 5).Value = "  +2.70%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.616.19"
$ws.Cells.Item(3, 5).Value = "  +0.90%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.04%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "571.18"
$ws.Cells.Item(5, 5).Value = "  -0.20%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "143.33"
$ws.Cells.Item(6, 5).Value = "  -0.54%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.26%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.20%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "2.639.54"
$ws.Cells.Item(9, 5).Value = "  +1.33%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "6.70"
$ws.Cells.Item(10, 5).Value = "  +0.38%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +1.90%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.08%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "0.376"
$ws.Cells.Item(13, 5).Value = "  +8.20%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "3.078.14"
$ws.Cells.Item(14, 5).Value = "  +0.98%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "60.981.14"
$ws.Cells.Item(15, 5).Value = "  +2.57%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "23.58"
$ws.Cells.Item(16, 5).Value = "  +4.23%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +1.89%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.626.02"
$ws.Cells.Item(18, 5).Value = "  +1.20%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "4.70"
$ws.Cells.Item(19, 5).Value = "  +3.24%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "Chainlink"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Cells.Item(20, 4) "11.05"
$ws.Cells.Item(20, 5).Value = "  +7.58%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "BitcoinCash"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Cells.Item(21, 4) "350.57"
$ws.Cells.Item(21, 5).Value = "  +3.66%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "7.08"
$ws.Cells.Item(22, 5).Value = "  +13.81%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.29%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "0.516"
$ws.Cells.Item(24, 5).Value = "  +13.29%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "63.63"
$ws.Cells.Item(25, 5).Value = "  -1.51%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "0.999"
$ws.Cells.Item(26, 5).Value = "  +0.33%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -0.06%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "7.75"
$ws.Cells.Item(28, 5).Value = "  +6.05%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "0.0₃0796"
$ws.Cells.Item(29, 5).Value = "  +1.45%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "1.85"
$ws.Cells.Item(30, 5).Value = "  +9.56%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.11%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "6.28"
$ws.Cells.Item(32, 5).Value = "  +2.96%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "161.60"
$ws.Cells.Item(33, 5).Value = "  +1.35%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) "19.49"
$ws.Cells.Item(34, 5).Value = "  +2.03%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "4.24"
$ws.Cells.Item(35, 5).Value = "  +4.56%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "0.965"
$ws.Cells.Item(36, 5).Value = "  +8.74%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +4.63%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "1.60"
$ws.Cells.Item(38, 5).Value = "  +6.53%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "37.69"
$ws.Cells.Item(39, 5).Value = "  +1.47%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.859"
$ws.Cells.Item(40, 5).Value = "  -2.25%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +3.15%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "303.89"
$ws.Cells.Item(42, 5).Value = "  +2.71%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "141.11"
$ws.Cells.Item(43, 5).Value = "  +13.59%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "0.993"
$ws.Cells.Item(44, 5).Value = "  -0.52%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "0.0986"
$ws.Cells.Item(45, 5).Value = "  +0.62%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  +1.55%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "0.0549"
$ws.Cells.Item(47, 5).Value = "  +1.71%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +3.81%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "10.70"
$ws.Cells.Item(49, 5).Value = "  +0.62%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "19.53"
$ws.Cells.Item(50, 5).Value = "  +4.94%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +7.16%  "
